$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: bxywxh3i.zso is different with handoff file name: c1761016-b1c4-4744-b8d7-9d90bf887f60.b1184278c36cbee1ed94481afb2799053227b076.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: bxywxh3i.zso is different with handoff file name: c1761016-b1c4-4744-b8d7-9d90bf887f60.b1184278c36cbee1ed94481afb2799053227b076.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
